# Students are not allowed to choose two projects that have the same
# supervisor. This reflects the roster after that validation was added:
# four students' project allocations were reassigned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Assunta Feil -> Cloyd Goldner, now has Applied/Accepted = 1
$ws.Range("A2").Value = "Cloyd Goldner"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3: Al Greenholt -> Gregg Grant (counts unchanged)
$ws.Range("A3").Value = "Gregg Grant"

# Row 4: Keven Mueller -> Destinee Orn, now has Applied/Accepted = 0
$ws.Range("A4").Value = "Destinee Orn"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5: Katharina Shields -> Alva Purdy (counts unchanged)
$ws.Range("A5").Value = "Alva Purdy"

# Recompute the "Name" column's best-fit width now that the longest name
# is shorter than before.
$ws.Columns.Item(1).AutoFit()
